# Invalid login Test Case
# Adds a new "InvalidLogin" worksheet (after "ValidLogin") containing a
# UserName/Password header row plus a sample invalid credential pair, and
# updates the selection/active-sheet state to match.

$wb = $excel.ActiveWorkbook

# --- ValidLogin (sheet 1): change selection to the A1:B2 range -------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A1:B2").Select() | Out-Null

# --- Add the new InvalidLogin sheet, placed right after ValidLogin ---------
$wsInvalid = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$wsInvalid.Name = "InvalidLogin"

$wsInvalid.Range("A1").Value = "UserName"
$wsInvalid.Range("B1").Value = "Password"
$wsInvalid.Range("A2").Value = "Bhanu"
$wsInvalid.Range("B2").Value = "Damager"

# Match the zoom level / selection captured in the new sheet's view.
$excel.ActiveWindow.Zoom = 250
$wsInvalid.Range("B3").Select() | Out-Null
